# Generate Report for Handback
# Update the localization status workbook: mark the handback as failed,
# record the error detail message and widen the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"
$zhMsg = "Handback file name: up4ic4wh.ma4 is different with handoff file name: 51284213-18c5-48ee-aca6-a34fa0ac1473.8ce6b9cd1b6a44a2338958ce56e369e547df6f41.zh-cn."
$deMsg = "Handback file name: up4ic4wh.ma4 is different with handoff file name: 51284213-18c5-48ee-aca6-a34fa0ac1473.8ce6b9cd1b6a44a2338958ce56e369e547df6f41.de-de."

# Overview sheet: the per-language status summary cells for the second
# data row ("Ready for handoff" everywhere this text occurs in the workbook).
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: update status and error detail for the second data row (row 3)
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = $zhMsg
# ColumnWidth round-trips through a pixel-based rounding, so 39.17 is the
# value that stores back to a character width of exactly 40 (matching the
# other width="40" columns already on this sheet).
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# de-de sheet: update status and error detail for the second data row (row 3)
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = $deMsg
$dede.Columns.Item(16).ColumnWidth = 39.17
